# Updated cryptos list on Mon Aug 12 06:58:22 UTC 2024 with GitHub Actions
#
# Refreshes the Price (column D) and Volume(1h) (column E) figures on the
# active worksheet for the rows whose crypto quotes moved since the last
# snapshot. Rows/columns not listed here are left untouched.
#
# Price values are stored as text in this sheet (e.g. "58.324.80" uses a
# dotted thousands style, and some, like "1.00"/"35.90", must keep their
# trailing zero). Assigning a plain numeric-looking string to a Range's
# Value lets Excel auto-convert it to a real number, which would silently
# normalise "1.00" -> 1 and "35.90" -> 35.9. To prevent that, those cells
# are briefly switched to a text number format before the write and then
# restored to the Normal style afterwards so no extra formatting lingers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.324.80"
$ws.Range("E2").Value = "  -4.46%  "

$ws.Range("D3").Value = "2.541.53"
$ws.Range("E3").Value = "  -3.92%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "505.64"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.98%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.42"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -7.83%  "

$ws.Range("E7").Value = "  +0.06%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.563"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.83%  "

$ws.Range("D9").Value = "2.544.40"
$ws.Range("E9").Value = "  -4.24%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.09"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -8.21%  "

$ws.Range("E11").Value = "  -7.14%  "

$ws.Range("E12").Value = "  -5.78%  "

$ws.Range("E13").Value = "  -0.76%  "

$ws.Range("D14").Value = "2.984.29"
$ws.Range("E14").Value = "  -4.05%  "

$ws.Range("D15").Value = "58.300.08"
$ws.Range("E15").Value = "  -4.48%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.65"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -6.41%  "

$ws.Range("E17").Value = "  -6.47%  "

$ws.Range("D18").Value = "2.545.20"
$ws.Range("E18").Value = "  -4.04%  "

$ws.Range("E19").Value = "  -5.51%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "341.68"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.70%  "

$ws.Range("E21").Value = "  -5.92%  "

$ws.Range("E22").Value = "  -0.05%  "

$ws.Range("E23").Value = "  -5.06%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "60.47"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.99%  "

$ws.Range("E25").Value = "  -5.13%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.11%  "

$ws.Range("D27").Value = "2.656.43"
$ws.Range("E27").Value = "  -3.95%  "

$ws.Range("E28").Value = "  -6.38%  "

$ws.Range("E29").Value = "  -9.12%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.92"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -6.34%  "

$ws.Range("E31").Value = "  -0.03%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "149.32"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.46%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.46"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.53%  "

$ws.Range("E34").Value = "  -6.23%  "

$ws.Range("E35").Value = "  -6.25%  "

$ws.Range("E36").Value = "  -6.86%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.899"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.51%  "

$ws.Range("E38").Value = "  -8.36%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "35.90"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.61%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.816"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -11.04%  "

$ws.Range("E41").Value = "  -7.87%  "

$ws.Range("E42").Value = "  -8.24%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "280.66"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -8.99%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.997"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.01%  "

$ws.Range("E45").Value = "  -3.13%  "

$ws.Range("E46").Value = "  -7.28%  "

$ws.Range("E47").Value = "  -5.77%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "18.58"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -6.51%  "

$ws.Range("E50").Value = "  -6.13%  "

$ws.Range("E51").Value = "  -9.57%  "
